# Refactor: flip Yes/No flags on the Interface sheet's E/F columns
# (mirrors a re-run of the test-data generator against the upgraded
# selenium libraries) and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interface")

# Column E (rows 3-9, 11): Yes -> No
$ws.Range("E3").Value = "No"
$ws.Range("E4").Value = "No"
$ws.Range("E5").Value = "No"
$ws.Range("E6").Value = "No"
$ws.Range("E7").Value = "No"
$ws.Range("E8").Value = "No"
$ws.Range("E9").Value = "No"
$ws.Range("E11").Value = "No"

# Column F: row 9 Yes -> No, row 10 No -> Yes
$ws.Range("F9").Value = "No"
$ws.Range("F10").Value = "Yes"

# Selection moves from I10 to E11 on the Interface sheet
$ws.Range("E11").Select()
